# Se elimina el EC anterior (ordenado por periodo ascendente, alternando
# trabajador) y se agrega el nuevo bloque de datos (agrupado por
# trabajador, con periodos en orden descendente), segun la base de
# datos actualizada.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Datos nuevos: Tipo Doc (B, sin cambios), N Doc Trabajador (C),
# Nombre Trabajador (D), Periodo Mora (E), Valor Mora (F)
$rows = @(
    @{ Row = 16; Doc = "33297383";  Nombre = "MARIA EUGENIA ROCHA PABUENA";  Periodo = "1903"; Valor = 26041 },
    @{ Row = 17; Doc = "33297383";  Nombre = "MARIA EUGENIA ROCHA PABUENA";  Periodo = "1902"; Valor = 31249 },
    @{ Row = 18; Doc = "33297383";  Nombre = "MARIA EUGENIA ROCHA PABUENA";  Periodo = "1901"; Valor = 31249 },
    @{ Row = 19; Doc = "33297383";  Nombre = "MARIA EUGENIA ROCHA PABUENA";  Periodo = "1812"; Valor = 31249 },
    @{ Row = 20; Doc = "33297383";  Nombre = "MARIA EUGENIA ROCHA PABUENA";  Periodo = "1811"; Valor = 31249 },
    @{ Row = 21; Doc = "45528190";  Nombre = "MARTHA LUCIA ROCHA PABUENA";   Periodo = "1903"; Valor = 26041 },
    @{ Row = 22; Doc = "45528190";  Nombre = "MARTHA LUCIA ROCHA PABUENA";   Periodo = "1902"; Valor = 31249 },
    @{ Row = 23; Doc = "45528190";  Nombre = "MARTHA LUCIA ROCHA PABUENA";   Periodo = "1901"; Valor = 31249 },
    @{ Row = 24; Doc = "45528190";  Nombre = "MARTHA LUCIA ROCHA PABUENA";   Periodo = "1812"; Valor = 31249 },
    @{ Row = 25; Doc = "45528190";  Nombre = "MARTHA LUCIA ROCHA PABUENA";   Periodo = "1811"; Valor = 31249 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("C$n").Value = $r.Doc
    $ws.Range("D$n").Value = $r.Nombre
    $ws.Range("E$n").Value = $r.Periodo
    $ws.Range("F$n").Value = $r.Valor
}
